# Insert a new weekly price row for "Agrícola del Norte S.A. de Arica - Albahaca"
# above the current row 49, shifting the existing rows 49:56 down to 50:57.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 49 downward (formatting of row 49 is carried along automatically).
$ws.Rows(49).Insert()

# Populate the newly inserted row 49 with the new weekly record.
$ws.Range("A49").Value = 1
$ws.Range("B49").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C49").Value = "Arica y Parinacota"
$ws.Range("D49").Value = 44995
$ws.Range("E49").Value = 15
$ws.Range("F49").Value = 100112052
$ws.Range("G49").Value = "Albahaca"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 500
$ws.Range("K49").Value = 1000
$ws.Range("L49").Value = 1000
$ws.Range("M49").Value = 1000
$ws.Range("N49").Value = "$/paquete"
$ws.Range("O49").Value = "Región de Arica y Parinacota"
$ws.Range("P49").Value = 1000
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = "Hortaliza"
